# Update "In Class Demonstration" worksheet with new FFR column values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (FFR) regression results being updated/added
$ws.Range("C2").Value = "-0.031***"
$ws.Range("C3").Value = "0.378***"

# Force this cell to remain text (like "-0.018" before it), not a number
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "0.086"

$ws.Range("C5").Value = 0.43
